$wb = $excel.ActiveWorkbook

# --- Sheet "ATS Accuracy": row 5 (🔥🔥 category) ---
$ws1 = $wb.Worksheets.Item("ATS Accuracy")
$ws1.Range("C5").Value = 6
$ws1.Range("D5").Value = 7
$ws1.Range("E5").Value = 85.7

# --- Sheet "Total Accuracy": row 5 (🔥🔥 category) ---
$ws2 = $wb.Worksheets.Item("Total Accuracy")
$ws2.Range("B5").Value = 4
$ws2.Range("D5").Value = 9
$ws2.Range("E5").Value = 55.6
